$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert 4 new data rows before the current totals row (row 16) ---
# This pushes the old totals row (16) and footer row (17) down to 20 and 21.
$ws.Range("A16:N19").EntireRow.Insert()

# --- Copy formatting from the last existing data row (15) into the new rows (16-19) ---
$ws.Range("A15:N15").Copy()
$ws.Range("A16:N19").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Restore row heights to match the data-row pattern ---
$ws.Rows("16").RowHeight = 24.75
$ws.Rows("17").RowHeight = 25.5
$ws.Rows("18").RowHeight = 25.5
$ws.Rows("19").RowHeight = 24.75
$ws.Rows("20").RowHeight = 26.25
$ws.Rows("21").RowHeight = 16.5

# --- Re-create the merged cell layout for the new rows ---
$ws.Range("B16:G16").Merge()
$ws.Range("H16:K16").Merge()
$ws.Range("L16:M16").Merge()

$ws.Range("B17:G17").Merge()
$ws.Range("H17:K17").Merge()
$ws.Range("L17:M17").Merge()

$ws.Range("B18:G18").Merge()
$ws.Range("H18:K18").Merge()
$ws.Range("L18:M18").Merge()

$ws.Range("B19:G19").Merge()
$ws.Range("H19:K19").Merge()
$ws.Range("L19:M19").Merge()

# --- Update existing rows 13-15 (values changed / items re-sorted) ---
$ws.Range("B13").Value = "سرنجات 10 سم"
$ws.Range("H13").Value = "-2:0"
$ws.Range("L13").Value = 8
$ws.Range("N13").Value = "2:0"

$ws.Range("B14").Value = "سرنجات 3 سم"
$ws.Range("H14").Value = "-2:0"
$ws.Range("L14").Value = 4
$ws.Range("N14").Value = "2:0"

$ws.Range("B15").Value = "سرنجات 5 سم"
$ws.Range("H15").Value = "-1:0"
$ws.Range("L15").Value = 2
$ws.Range("N15").Value = "1:0"

# --- Fill the 4 new rows (16-19) ---
$ws.Range("A16").Value = 13
$ws.Range("B16").Value = "شفرات فينوس حريمي "
$ws.Range("H16").Value = "16:0"
$ws.Range("L16").Value = 40
$ws.Range("N16").Value = "2:0"

$ws.Range("A17").Value = 14
$ws.Range("B17").Value = "كالونا "
$ws.Range("H17").Value = "-1:0"
$ws.Range("L17").Value = 15
$ws.Range("N17").Value = "1:0"

$ws.Range("A18").Value = 15
$ws.Range("B18").Value = "كريم فيرند لافلي الصغير"
$ws.Range("H18").Value = "6:0"
$ws.Range("L18").Value = 20
$ws.Range("N18").Value = "1:0"

$ws.Range("A19").Value = 16
$ws.Range("B19").Value = "محلول ملح"
$ws.Range("H19").Value = "27:0"
$ws.Range("L19").Value = 48
$ws.Range("N19").Value = "2:0"

# --- Update the grand total (now in row 20) ---
$ws.Range("K20").Value = 550

Write-Host "done"
